# Fix sheet1 content per commit diff: insert a row at 13, correct misaligned
# content that shifted down, and populate the newly authored PT objectives,
# short-syllabus, full syllabus and bibliography text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $text, $col) {
    $rng = $ws.Range($addr)
    $rng.Font.Bold = $false
    $rng.Value = $text
    $rng.WrapText = $true
    $rng.VerticalAlignment = -4160
    if ($col -eq "C") {
        $rng.Font.Color = 255
    }
}

# Insert a new blank row at 13 (shifts old rows 13-23 down to 14-24),
# the newly blank A13 does not carry a label in the target layout.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()

# Row 10 (Objetivos:) B/C - was misaligned docente name, now the PT objectives text
Set-Cell "B10" "Apresentar os fundamentos da Contabilidade Geral, de modo que os alunos possam compreender os conceitos e princípios básicos de Contabilidade. Capacitar na produção e análise de documentos e de registros de transações contábeis. Compreender os sistemas de custeio e sua aplicação. Capacitar a projetar e implantar sistemas de custeio. Capacitar a analisar custos em relação a estratégia de mercado e de produção de uma organização." "B"

# 
Set-Cell "C10" "Apresentar os fundamentos da Contabilidade Geral, de modo que os alunos possam compreender os conceitos e princípios básicos de Contabilidade. Capacitar na produção e análise de documentos e de registros de transações contábeis. Compreender os sistemas de custeio e sua aplicação. Capacitar a projetar e implantar sistemas de custeio. Capacitar a analisar custos em relação a estratégia de mercado e de produção de uma organização." "C"

# Row 13 (blank label) B/C - docente responsible, shifted down from row 10's old slot
Set-Cell "B13" "11079086 - Herlandí de Souza Andrade" "B"

# 
Set-Cell "C13" "11079086 - Herlandí de Souza Andrade" "C"

# Row 14 (Programa resumido:) B/C - short syllabus PT (was 'Semestral')
Set-Cell "B14" "1  Fundamentos de Contabilidade. 2 - Sistemas de Custeio. 3 - Projeto de Sistemas de Custeio. 4 - Análise de Custos." "B"

# 
Set-Cell "C14" "1  Fundamentos de Contabilidade. 2 - Sistemas de Custeio. 3 - Projeto de Sistemas de Custeio. 4 - Análise de Custos." "C"

# Row 16 (Programa:) B/C - full syllabus PT content (was stray date '01/01/2021')
Set-Cell "B16" "1. FUNDAMENTOS DE CONTABILIDADEConceito de contabilidade. Contabilidade e Engenharia de Produção. Balanço: ativo, passivo, patrimônio líquido. Procedimentos contábeis básicos. Variações da situação líquida; despesa, receita. Regimes de competência e caixa. Receitas e despesas diferidas. Fatos Contábeis. Operações com mercadorias. Inventário. Ativo Imobilizado e Amortização. Demonstrativo de origens e aplicações. Análise de Balanço e de resultados.2. SISTEMAS DE CUSTEIOCustos diretos e indiretos, fixos e variáveis. Acumulação de custos, classificação, fatores de custo. Custeio por Absorção. Produção por ordem, contínua, conjunta. Custeio Direto: margem de contribuição. Custeio ABC.3. PROJETO DE SISTEMAS DE CUSTEIOProdutos e Departamentos. Sistemas de produção e sistemas de custeio. O problema da inflação.4. ANÁLISE DE CUSTOSCusto-Volume-Lucro. Contribuição marginal. Análise de variações. Equação de produtividade global. Alavancagem operacional. TIR e lucratividade." "B"

# 
Set-Cell "C16" "1. FUNDAMENTOS DE CONTABILIDADEConceito de contabilidade. Contabilidade e Engenharia de Produção. Balanço: ativo, passivo, patrimônio líquido. Procedimentos contábeis básicos. Variações da situação líquida; despesa, receita. Regimes de competência e caixa. Receitas e despesas diferidas. Fatos Contábeis. Operações com mercadorias. Inventário. Ativo Imobilizado e Amortização. Demonstrativo de origens e aplicações. Análise de Balanço e de resultados.2. SISTEMAS DE CUSTEIOCustos diretos e indiretos, fixos e variáveis. Acumulação de custos, classificação, fatores de custo. Custeio por Absorção. Produção por ordem, contínua, conjunta. Custeio Direto: margem de contribuição. Custeio ABC.3. PROJETO DE SISTEMAS DE CUSTEIOProdutos e Departamentos. Sistemas de produção e sistemas de custeio. O problema da inflação.4. ANÁLISE DE CUSTOSCusto-Volume-Lucro. Contribuição marginal. Análise de variações. Equação de produtividade global. Alavancagem operacional. TIR e lucratividade." "C"

# Row 19 (Metodo:) B/C - method text, shifted down from row 18's old slot
Set-Cell "B19" "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras." "B"

# 
Set-Cell "C19" "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras." "C"

# Row 20 (Criterio:) B/C - criterion text, shifted down from row 19's old slot
Set-Cell "B20" "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas." "B"

# 
Set-Cell "C20" "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas." "C"

# Row 21 (Norma de recuperacao:) B/C - recovery-exam formula, shifted from row 20
Set-Cell "B21" "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação." "B"

# 
Set-Cell "C21" "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação." "C"

# Row 22 (Bibliografia:) B/C - full bibliography text (was stray recovery formula)
Set-Cell "B22" "IUDICIBUS, S.; MARION, J. C. Curso de Contabilidade para não Contadores. 8 ed. São Paulo: Atlas, 2018.MARTINS, E. Contabilidade de Custos. 11 ed. São Paulo: Atlas, 2018.SANTOS, J. J. Manual de Contabilidade e Análise de Custos. 7 ed. São Paulo: Atlas, 2017.SILVA, R. N. S.; LINS, L. S. Gestão de Custos - Contabilidade, Controle e Análise. 4 ed. São Paulo: Atlas, 2017.DUTRA, R. G. Custos: Uma Abordagem Prática. 8 ed. São Paulo: Atlas, 2017.MEGLIORINI, E.; BUENO, A. S. Contabilidade para cursos de Engenharia. São Paulo: Atlas, 2014.MORANTE, A. S. Análise das Demonstrações Financeiras. 2 ed. São Paulo: Atlas, 2009.Bruni, Adriano L., Fama, Rubens. Gestão de custos e formação de preços. 5. Ed., São Paulo: Atlas, 2008.Hansen, Don R., Mowen, Maryanne M., Gestão de custos – contabilidade e controle. São Paulo: Pioneira Thomson, 2001Horngren, Charles T., Datar, Srikant M., Foster, George. Contabilidade de custos (vol. 1 e 2). 11. Ed. São Paulo: Pearson, 2004.Maher, Michael. Contabilidade de custos – criando valor para a administração. 5. ed., São Paulo: Atlas, 2001.Martins, Eliseu. Contabilidade de Custos (livro de exercícios). 9. Ed., São Paulo: Atlas, 2006.Martins, Eliseu. Contabilidade de Custos (livro texto). 9. Ed., São Paulo: Atlas, 2003." "B"

# 
Set-Cell "C22" "IUDICIBUS, S.; MARION, J. C. Curso de Contabilidade para não Contadores. 8 ed. São Paulo: Atlas, 2018.MARTINS, E. Contabilidade de Custos. 11 ed. São Paulo: Atlas, 2018.SANTOS, J. J. Manual de Contabilidade e Análise de Custos. 7 ed. São Paulo: Atlas, 2017.SILVA, R. N. S.; LINS, L. S. Gestão de Custos - Contabilidade, Controle e Análise. 4 ed. São Paulo: Atlas, 2017.DUTRA, R. G. Custos: Uma Abordagem Prática. 8 ed. São Paulo: Atlas, 2017.MEGLIORINI, E.; BUENO, A. S. Contabilidade para cursos de Engenharia. São Paulo: Atlas, 2014.MORANTE, A. S. Análise das Demonstrações Financeiras. 2 ed. São Paulo: Atlas, 2009.Bruni, Adriano L., Fama, Rubens. Gestão de custos e formação de preços. 5. Ed., São Paulo: Atlas, 2008.Hansen, Don R., Mowen, Maryanne M., Gestão de custos – contabilidade e controle. São Paulo: Pioneira Thomson, 2001Horngren, Charles T., Datar, Srikant M., Foster, George. Contabilidade de custos (vol. 1 e 2). 11. Ed. São Paulo: Pearson, 2004.Maher, Michael. Contabilidade de custos – criando valor para a administração. 5. ed., São Paulo: Atlas, 2001.Martins, Eliseu. Contabilidade de Custos (livro de exercícios). 9. Ed., São Paulo: Atlas, 2006.Martins, Eliseu. Contabilidade de Custos (livro texto). 9. Ed., São Paulo: Atlas, 2003." "C"
